# Update countries & provincias Spain
# Applies the 26-May-2020 14:35 data refresh to the "Pais" sheet:
#  - updates the "last updated" timestamp
#  - refreshes case counters for several countries
#  - a handful of countries swap rank (and therefore row position) because
#    their updated "Casos totales" (col B) changed their sort order

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Timestamp banner (row 1) --------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 26 de Mayo de 2020 a las 14:35"

# ---- Simple in-place numeric refreshes (no reordering) -------------------

# India (row 13)
$ws.Range("B13").Value = 146443
$ws.Range("C13").Value = 1493
$ws.Range("E13").Value = 81105

# Croacia (row 87)
$ws.Range("D87").Value = 2046
$ws.Range("E87").Value = 97
$ws.Range("G87").Value = 1
$ws.Range("H87").Value = 101

# Libano (row 106)
$ws.Range("B106").Value = 1140
$ws.Range("C106").Value = 21
$ws.Range("D106").Value = 689
$ws.Range("E106").Value = 425

# ---- Guinea / Uzbekistan swap rank (rows 76-77) ---------------------------
# Uzbekistan's updated totals now outrank Guinea, so Uzbekistan moves to 76.
$ws.Range("A76").Value = "Uzbekistan"
$ws.Range("B76").Value = 3281
$ws.Range("C76").Value = 92
$ws.Range("D76").Value = 2624
$ws.Range("E76").Value = 644
$ws.Range("F76").Value = 0
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 13

$ws.Range("A77").Value = "Guinea"
$ws.Range("B77").Value = 3275
$ws.Range("C77").Value = 0
$ws.Range("D77").Value = 1673
$ws.Range("E77").Value = 1582
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 20

# ---- Benin / Birmania / Martinica reshuffle (rows 156-158) ----------------
# Benin's updated totals push it above both Birmania and Martinica.
$ws.Range("A156").Value = "Benin"
$ws.Range("B156").Value = 208
$ws.Range("C156").Value = 17
$ws.Range("D156").Value = 118
$ws.Range("E156").Value = 87
$ws.Range("F156").Value = 0
$ws.Range("G156").Value = 0
$ws.Range("H156").Value = 3

$ws.Range("A157").Value = "Birmania"
$ws.Range("B157").Value = 203
$ws.Range("C157").Value = 0
$ws.Range("D157").Value = 123
$ws.Range("E157").Value = 74
$ws.Range("F157").Value = 0
$ws.Range("G157").Value = 0
$ws.Range("H157").Value = 6

$ws.Range("A158").Value = "Martinica"
$ws.Range("B158").Value = 197
$ws.Range("C158").Value = 0
$ws.Range("D158").Value = 91
$ws.Range("E158").Value = 92
$ws.Range("F158").Value = 0
$ws.Range("G158").Value = 0
$ws.Range("H158").Value = 14

# ---- Belice / Santa Lucia swap rank (rows 200-201) -------------------------
$ws.Range("A200").Value = "Santa Lucia"
$ws.Range("B200").Value = 18
$ws.Range("C200").Value = 0
$ws.Range("D200").Value = 18
$ws.Range("E200").Value = 0
$ws.Range("F200").Value = 0
$ws.Range("G200").Value = 0
$ws.Range("H200").Value = 0

$ws.Range("A201").Value = "Belice"
$ws.Range("B201").Value = 18
$ws.Range("C201").Value = 0
$ws.Range("D201").Value = 16
$ws.Range("E201").Value = 0
$ws.Range("F201").Value = 0
$ws.Range("G201").Value = 0
$ws.Range("H201").Value = 2

# ---- Groenlandia / Islas Turcas y Caicos swap rank (rows 207-208) ---------
$ws.Range("A207").Value = "Islas Turcas y Caicos"
$ws.Range("B207").Value = 12
$ws.Range("C207").Value = 0
$ws.Range("D207").Value = 10
$ws.Range("E207").Value = 1
$ws.Range("F207").Value = 0
$ws.Range("G207").Value = 0
$ws.Range("H207").Value = 1

$ws.Range("A208").Value = "Groenlandia"
$ws.Range("B208").Value = 12
$ws.Range("C208").Value = 0
$ws.Range("D208").Value = 11
$ws.Range("E208").Value = 1
$ws.Range("F208").Value = 0
$ws.Range("G208").Value = 0
$ws.Range("H208").Value = 0
